$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.106.04'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '3.885.56'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '484.40'
$ws.Range('E5').Value = '  +1.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.76'
$ws.Range('E6').Value = '  -1.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  +1.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.997'
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.738'
$ws.Range('E9').Value = '  +1.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.177'
$ws.Range('E10').Value = '  +6.87%  '
$ws.Range('E11').Value = '  +1.55%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.95'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.60'
$ws.Range('E13').Value = '  +3.26%  '
$ws.Range('D14').Value = '4.508.45'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('D15').Value = '3.889.39'
$ws.Range('E15').Value = '  -1.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.32'
$ws.Range('E16').Value = '  -2.09%  '
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.26'
$ws.Range('E18').Value = '  +2.50%  '
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('D20').Value = '68.096.40'
$ws.Range('E20').Value = '  +0.55%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '429.97'
$ws.Range('E21').Value = '  -0.64%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.56'
$ws.Range('E22').Value = '  +6.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.83'
$ws.Range('E23').Value = '  +3.48%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.21'
$ws.Range('E24').Value = '  +2.52%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.61'
$ws.Range('E25').Value = '  +13.91%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.69'
$ws.Range('E26').Value = '  +4.72%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.05'
$ws.Range('E27').Value = '  +9.87%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.46'
$ws.Range('E28').Value = '  -1.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.68'
$ws.Range('E29').Value = '  -2.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '719.59'
$ws.Range('E30').Value = '  -1.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.66'
$ws.Range('E31').Value = '  +3.28%  '
$ws.Range('E32').Value = '  +1.28%  '
$ws.Range('E33').Value = '  +3.67%  '
$ws.Range('D34').Value = '0.0₃0886'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '41.50'
$ws.Range('E35').Value = '  -1.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.09'
$ws.Range('E36').Value = '  +13.47%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '61.24'
$ws.Range('E37').Value = '  +4.29%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.398'
$ws.Range('E38').Value = '  +17.30%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.146'
$ws.Range('E39').Value = '  -3.31%  '
$ws.Range('B40').Value = 'Dai'
$ws.Range('C40').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.999'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.02'
$ws.Range('E41').Value = '  +8.10%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0496'
$ws.Range('E42').Value = '  +6.18%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.11'
$ws.Range('E43').Value = '  +4.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.98'
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.143'
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.37'
$ws.Range('E46').Value = '  +6.49%  '
$ws.Range('E47').Value = '  +0.37%  '
$ws.Range('E48').Value = '  -0.93%  '
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '144.70'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.81'
$ws.Range('E51').Value = '  -1.52%  '
